$p = $ppt.ActivePresentation

# --- 1. Remove the "XmlAdaptedTag" rectangle (shape id 73) and its
#        connecting elbow connector (shape id 70) from the storage
#        class diagram on slide 8. ---
$s8 = $p.Slides.Item(8)

for ($i = $s8.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s8.Shapes.Item($i)
    if ($shp.Id -eq 70 -or $shp.Id -eq 73) {
        $shp.Delete()
    }
}

# --- 2. Refresh the cached "datetimeFigureOut" footer text (slide
#        master + every layout) from 10/14/16 to 10/15/2016. ---
function Update-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "10/14/16") {
            $shp.TextFrame.TextRange.Text = "10/15/2016"
        }
    }
}

$m = $p.SlideMaster
Update-DateText $m.Shapes
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    Update-DateText $m.CustomLayouts.Item($li).Shapes
}
